# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de sheets (row 2 entries) to reflect a
# freshly-generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 17:15:48"
$wsZhCn.Range("H2").Value = "2016-03-22 17:16:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 17:15:53"
$wsDeDe.Range("H2").Value = "2016-03-22 17:16:27"
